$d = $word.ActiveDocument

# --- 1) Bold + underline the "volumineux" duplicate-folders heading ------
$headingRange = $d.Content
$headingFound = $headingRange.Find.Execute("Top 5 des éléments présents plusieurs fois les plus volumineux", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($headingFound) {
    $paraCount = $d.Paragraphs.Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -eq $headingRange.Start) {
            # Bold the whole paragraph, including its paragraph mark, so the
            # paragraph-mark run properties (w:pPr/w:rPr) pick up <w:b/> too.
            $p.Range.Bold = 1
            break
        }
    }
    # Underline only the visible text run (not the paragraph mark).
    $headingRange2 = $d.Content
    $headingRange2.Find.Execute("Top 5 des éléments présents plusieurs fois les plus volumineux", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $headingRange2.Underline = 1
}

# --- 2) Text replacements: "files" -> "folders" ---------------------------
$d.Content.Find.Execute("{#biggestDuplicateFiles}", $true, $false, $false, $false, $false, $true, 1, $false, "{#biggestDuplicateFolders}", 2) | Out-Null
$d.Content.Find.Execute("Chemin du fichier ayant la date de modification la plus ancienne", $true, $false, $false, $false, $false, $true, 1, $false, "Chemin du dossier ayant la date de modification la plus ancienne", 2) | Out-Null
$d.Content.Find.Execute("{/biggestDuplicateFiles}", $true, $false, $false, $false, $false, $true, 1, $false, "{/biggestDuplicateFolders}", 2) | Out-Null

Write-Output "done"
